$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.471327662467957
$ws.Range("B1").Value = 1.583739995956421
$ws.Range("C1").Value = 1.776378273963928
$ws.Range("D1").Value = 2.732922315597534
$ws.Range("E1").Value = 3.515083312988281
